$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values must remain plain text, matching the
# original inline-string cells, so force Text number format before assigning).
$updates = @{
    "D2" = "258.48"
    "E2" = "0.26%"
    "E3" = "-0.37%"
    "D4" = "4.614"
    "E4" = "-6.06%"
    "D5" = "0.05909"
    "E5" = "-0.68%"
    "D6" = "6.642"
    "E6" = "-0.75%"
    "D7" = "0.8534"
    "E7" = "-1.90%"
    "D8" = "0.9467"
    "E8" = "-4.14%"
    "D9" = "0.1403"
    "E9" = "-0.42%"
    "D10" = "0.05237"
    "E10" = "47.35%"
    "D11" = "0.07097"
    "E11" = "-0.97%"
    "D12" = "0.03106"
    "E12" = "-1.36%"
    "D13" = "0.09145"
    "E13" = "-1.18%"
    "D14" = "0.001524"
    "E14" = "-1.94%"
    "D15" = "0.006116"
    "E15" = "1.75%"
    "D16" = "3.504"
    "E16" = "0.53%"
    "D17" = "3.190"
    "E17" = "-2.16%"
    "E18" = "-0.08%"
    "D19" = "0.0006048"
    "E19" = "-0.85%"
    "E20" = "-2.88%"
    "E21" = "-2.20%"
    "E22" = "8.18%"
    "D23" = "0.04258"
    "E23" = "-0.08%"
    "D24" = "0.001220"
    "E24" = "-0.16%"
    "D25" = "0.004296"
    "E25" = "-4.92%"
    "E26" = "-0.03%"
    "D27" = "0.0001936"
    "E27" = "29.86%"
    "D40" = "0.03826"
    "E40" = "-0.21%"
    "D41" = "0.006238"
    "E41" = "57.18%"
    "D42" = "0.1103"
    "E42" = "-0.14%"
    "D43" = "0.002199"
    "E43" = "-4.79%"
    "D44" = "0.01412"
    "E44" = "34.13%"
    "D45" = "0.00005408"
    "E45" = "-1.56%"
    "E46" = "-0.03%"
    "D47" = "0.05099"
    "E47" = "-53.26%"
    "D48" = "0.2524"
    "E48" = "11,771.46%"
    "D49" = "0.00002099"
    "E49" = "-0.03%"
    "D50" = "0.0001999"
    "E50" = "-0.03%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
